# C5-PowerPoint.pptx edit: 30 Jul 2020
#
# 1) The table on slide 6 gets a new (built-in) table style applied.
# 2) The deck's theme colour scheme is swapped from the custom "Integral"
#    palette to the stock Office palette (dk1/lt1 unchanged, the other ten
#    scheme colours change).

$p = $ppt.ActivePresentation

# --- 1. Re-style the table on slide 6 -------------------------------------
$slide6 = $p.Slides.Item(6)
for ($i = 1; $i -le $slide6.Shapes.Count; $i++) {
    $shp = $slide6.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{6596B210-A375-42B7-9E70-C797F27A5755}")
    }
}

# --- 2. Swap the theme colour scheme over to the stock Office palette -----
$tcs = $slide6.ThemeColorScheme
$tcs.Colors(1).RGB  = 0         # dk1      000000
$tcs.Colors(2).RGB  = 16777215  # lt1      FFFFFF
$tcs.Colors(3).RGB  = 6968388   # dk2      44546A
$tcs.Colors(4).RGB  = 15132391  # lt2      E7E6E6
$tcs.Colors(5).RGB  = 13998939  # accent1  5B9BD5
$tcs.Colors(6).RGB  = 3243501   # accent2  ED7D31
$tcs.Colors(7).RGB  = 10855845  # accent3  A5A5A5
$tcs.Colors(8).RGB  = 49407     # accent4  FFC000
$tcs.Colors(9).RGB  = 12874308  # accent5  4472C4
$tcs.Colors(10).RGB = 4697456   # accent6  70AD47
$tcs.Colors(11).RGB = 12673797  # hlink    0563C1
$tcs.Colors(12).RGB = 7491477   # folHlink 954F72
